$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "127.0.0.1"
$ws.Range("F3").Value = "127.0.0.1"
$ws.Range("F4").Value = "127.0.0.1"
$ws.Range("F5").Value = "127.0.0.1"
$ws.Range("F6").Value = "127.0.0.1"

$ws.Range("F14").Select()
